$d = $word.ActiveDocument
$tbl = $d.Tables.Item(2)

# Find the data row describing "Nguyễn Hữu Hòa" (STT = 6, nhiệm kỳ 2021 - 2023) so the
# whole table row can be removed, regardless of its exact position.
$targetIndex = -1
for ($i = 1; $i -le $tbl.Rows.Count; $i++) {
    $infoCell = $tbl.Rows.Item($i).Cells.Item(2)
    if ($infoCell.Range.Text -like "*Nguyễn Hữu Hòa*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ge 1) {
    # Remember the STT label of that row (e.g. "6") before the row disappears.
    $sttLabel = $tbl.Rows.Item($targetIndex).Cells.Item(1).Range.Text.Trim()

    # Delete the entire row; every following row shifts up to fill the gap.
    $tbl.Rows.Item($targetIndex).Delete()

    # The row that shifted into the deleted row's place still carries its own
    # (now one-too-high) STT number; relabel it with the vacated number so the
    # STT column stays sequential.
    $shiftedCell = $tbl.Rows.Item($targetIndex).Cells.Item(1)
    $shiftedCell.Range.Text = $sttLabel
}
